# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" column (D) on each language sheet so
# every row reflects the most recent handoff run for that language.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D2:D5").Value = "2016-02-29 13:51:30"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D2:D5").Value = "2016-02-29 13:51:45"
